# DONALD-4 all files for P1 and P2
# Fix a typo in the variable names "m_berufsab" -> "m_berufab" and
# "v_berufsab" -> "v_berufab" (dropping the extra "s" from "berufsab")
# on both the "Variables" sheet and the "Categories" sheet.

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCategories = $wb.Worksheets.Item("Categories")

# --- Variables sheet --------------------------------------------------
# Row 7: variable name "b_berufsab" -> "m_berufab"
$wsVariables.Range("B7").Value = "m_berufab"
# Row 8: variable name "v_berufsab" -> "v_berufab"
$wsVariables.Range("B8").Value = "v_berufab"

# --- Categories sheet ---------------------------------------------------
# Rows 14-21: the "variable" column referenced "m_berufsab" -> "m_berufab"
$wsCategories.Range("A14:A21").Value = "m_berufab"
# Rows 22-29: the "variable" column referenced "v_berufsab" -> "v_berufab"
$wsCategories.Range("A22:A29").Value = "v_berufab"
